# Regenerate save_data column G ("K") values: replace old Strike# values
# with new K values (recalculated "s_vals").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 5
    4  = 1
    5  = 2
    6  = 4
    7  = 3
    8  = 1
    9  = 5
    10 = 1
    11 = 6
    12 = 1
    13 = 2
    14 = 3
    15 = 4
    16 = 5
    17 = 5
    18 = 0
    19 = 8
    20 = 5
    21 = 4
    22 = 4
    23 = 2
    24 = 6
    25 = 4
    26 = 7
    27 = 6
    28 = 7
    29 = 7
    30 = 4
    31 = 4
    32 = 3
    33 = 2
    34 = 5
    35 = 2
    36 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
